$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 674
$ws.Range("I12").Value = 349
$ws.Range("K12").Value = 349
$ws.Range("M12").Value = -179
$ws.Range("H17").Value = 795.1458
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 795.1458
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2385.4374
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2721.4374
$ws.Range("H108").Value = 57000
$ws.Range("J108").Value = 57000
$ws.Range("L108").Value = 57000
$ws.Range("N108").Value = -64680
$ws.Range("H116").Value = 3732.5625
$ws.Range("I116").Value = 3677
$ws.Range("K116").Value = 3677
$ws.Range("M116").Value = -235
$ws.Range("H138").Value = 45457052
$ws.Range("J138").Value = 83335384
$ws.Range("L138").Value = 250006152
$ws.Range("N138").Value = -250016432
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 1361.1428
$ws.Range("I13").Value = 925.6
$ws.Range("J13").Value = 2450
$ws.Range("K13").Value = 925.6
$ws.Range("L13").Value = 2450
$ws.Range("M13").Value = -781.6
$ws.Range("N13").Value = -2738
$ws.Range("H45").Value = 5878.2144
$ws.Range("I45").Value = 3470.7144
$ws.Range("J45").Value = 8285.714
$ws.Range("K45").Value = 3470.7144
$ws.Range("L45").Value = 8285.714
$ws.Range("M45").Value = -3093.7144
$ws.Range("N45").Value = -9039.714
$ws.Range("H97").Value = 1460.6364
$ws.Range("I97").Value = 1373.2858
$ws.Range("J97").Value = 1613.5
$ws.Range("K97").Value = 1373.2858
$ws.Range("L97").Value = 1613.5
$ws.Range("M97").Value = -877.2858000000001
$ws.Range("N97").Value = -2605.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 8500
$ws.Range("I22").Value = 8000
$ws.Range("K22").Value = 8000
$ws.Range("M22").Value = -7827
$ws.Range("H99").Value = 4286.8213
$ws.Range("I99").Value = 2917.4736
$ws.Range("K99").Value = 2917.4736
$ws.Range("M99").Value = -1419.4736
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2828
$ws.Range("I2").Value = 5106
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 5106
$ws.Range("L2").Value = 550
$ws.Range("M2").Value = -4993
$ws.Range("N2").Value = -776
$ws.Range("H68").Value = 41283.75
$ws.Range("J68").Value = 41283.75
$ws.Range("L68").Value = 41283.75
$ws.Range("N68").Value = -42781.75
$ws.Range("H71").Value = 41283.75
$ws.Range("J71").Value = 41283.75
$ws.Range("L71").Value = 123851.25
$ws.Range("N71").Value = -131339.25
$ws.Range("H107").Value = 487.09677
$ws.Range("I107").Value = 414.2
$ws.Range("J107").Value = 619.63635
$ws.Range("K107").Value = 414.2
$ws.Range("L107").Value = 619.63635
$ws.Range("M107").Value = 1505.8
$ws.Range("N107").Value = -4459.63635
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 20518742
$ws.Range("I4").Value = 18834768
$ws.Range("J4").Value = 41400000
$ws.Range("K4").Value = 56504304
$ws.Range("L4").Value = 124200000
$ws.Range("M4").Value = -56504192
$ws.Range("N4").Value = -124200224
$ws.Range("H11").Value = 7692924
$ws.Range("I11").Value = 9091473
$ws.Range("J11").Value = 901
$ws.Range("K11").Value = 27274419
$ws.Range("L11").Value = 2703
$ws.Range("M11").Value = -27274279
$ws.Range("N11").Value = -2983
$ws.Range("H12").Value = 1253.3572
$ws.Range("J12").Value = 1387.25
$ws.Range("L12").Value = 4161.75
$ws.Range("N12").Value = -4507.75
$ws.Range("H128").Value = 523357.9
$ws.Range("I128").Value = 523357.9
$ws.Range("K128").Value = 1570073.7
$ws.Range("M128").Value = -1565093.7
$ws.Range("H131").Value = 23812754
$ws.Range("I131").Value = 62500852
$ws.Range("K131").Value = 187502556
$ws.Range("M131").Value = -187497516
$ws.Range("H138").Value = 1439.8334
$ws.Range("I138").Value = 1439.8334
$ws.Range("K138").Value = 4319.5002
$ws.Range("M138").Value = 820.4997999999996
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 203.33333
$ws.Range("I2").Value = 144.8
$ws.Range("K2").Value = 144.8
$ws.Range("M2").Value = -31.80000000000001
$ws.Range("H3").Value = 4003334.8
$ws.Range("I3").Value = 6000000
$ws.Range("J3").Value = 10004
$ws.Range("K3").Value = 6000000
$ws.Range("L3").Value = 10004
$ws.Range("M3").Value = -5999884
$ws.Range("N3").Value = -10236
$ws.Range("H132").Value = 4020.4211
$ws.Range("J132").Value = 5577.4
$ws.Range("L132").Value = 16732.2
$ws.Range("N132").Value = -21792.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1928.0834
$ws.Range("I16").Value = 1797.6666
$ws.Range("K16").Value = 1797.6666
$ws.Range("M16").Value = -1627.6666
$ws.Range("H20").Value = 26499.95
$ws.Range("J20").Value = 36666.582
$ws.Range("L20").Value = 36666.582
$ws.Range("N20").Value = -37118.582
$ws.Range("H22").Value = 3555.5
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 3555.5
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H40").Value = 3234
$ws.Range("I40").Value = 2683.6667
$ws.Range("K40").Value = 2683.6667
$ws.Range("M40").Value = -2547.6667
$ws.Range("H61").Value = 1461
$ws.Range("I61").Value = 1461
$ws.Range("K61").Value = 1461
$ws.Range("M61").Value = -1259
$ws.Range("H68").Value = 1686.4546
$ws.Range("I68").Value = 1692.2858
$ws.Range("J68").Value = 1676.25
$ws.Range("K68").Value = 1692.2858
$ws.Range("L68").Value = 1676.25
$ws.Range("M68").Value = -943.2858000000001
$ws.Range("N68").Value = -3174.25
$ws.Range("H71").Value = 1686.4546
$ws.Range("I71").Value = 1692.2858
$ws.Range("J71").Value = 1676.25
$ws.Range("K71").Value = 8461.429
$ws.Range("L71").Value = 8381.25
$ws.Range("M71").Value = -4717.429
$ws.Range("N71").Value = -15869.25
$ws.Range("H113").Value = 1461
$ws.Range("I113").Value = 1461
$ws.Range("K113").Value = 1461
$ws.Range("M113").Value = 709
$ws.Range("H122").Value = 3370.5
$ws.Range("I122").Value = 3370.5
$ws.Range("K122").Value = 10111.5
$ws.Range("M122").Value = -7661.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 11931.214
$ws.Range("J41").Value = 12156.692
$ws.Range("L41").Value = 12156.692
$ws.Range("N41").Value = -12936.692
$ws.Range("H70").Value = 39576.668
$ws.Range("I70").Value = 35000
$ws.Range("K70").Value = 35000
$ws.Range("M70").Value = -34685
$ws.Range("H73").Value = 39576.668
$ws.Range("I73").Value = 35000
$ws.Range("K73").Value = 35000
$ws.Range("M73").Value = -33908
$ws.Range("H107").Value = 724.4706
$ws.Range("I107").Value = 558.4286
$ws.Range("K107").Value = 1675.2858
$ws.Range("M107").Value = 244.7142000000001
$ws.Range("H132").Value = 3866.348
$ws.Range("J132").Value = 3149.75
$ws.Range("L132").Value = 9449.25
$ws.Range("N132").Value = -14509.25
